$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "mannwhitneyu"
$ws.Range("E2").Value = 252
$ws.Range("F2").Value = 0.1895093944814901
$ws.Range("G2").Value = 0.2222222222222222
$ws.Range("H2").Value = 0.1864323930634537
$ws.Range("I2").Value = "Dunn"
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = -5.87777777777778

# Row 3
$ws.Range("D3").Value = "mannwhitneyu"
$ws.Range("E3").Value = 194.5
$ws.Range("F3").Value = 0.01792725432044502
$ws.Range("G3").Value = 0.3996913580246914
$ws.Range("H3").Value = 0.01748761387370319
$ws.Range("I3").Value = "Dunn"
$ws.Range("J3").Value = $true
$ws.Range("K3").Value = -2.585444444444446

# Row 4
$ws.Range("D4").Value = "mannwhitneyu"
$ws.Range("E4").Value = 186
$ws.Range("F4").Value = 0.01163445170501903
$ws.Range("G4").Value = 0.4259259259259259
$ws.Range("H4").Value = 0.01133439323885829
$ws.Range("I4").Value = "Dunn"
$ws.Range("J4").Value = $true
$ws.Range("K4").Value = -9.089416666666665

# Row 5
$ws.Range("D5").Value = "mannwhitneyu"
$ws.Range("E5").Value = 155
$ws.Range("F5").Value = 0.001988977628350087
$ws.Range("G5").Value = 0.5216049382716049
$ws.Range("H5").Value = 0.00192837068649817
$ws.Range("I5").Value = "Dunn"
$ws.Range("J5").Value = $true
$ws.Range("K5").Value = -1380.625666666665

# Row 9
$ws.Range("D9").Value = "mannwhitneyu"
$ws.Range("E9").Value = 243
$ws.Range("F9").Value = 0.1396266720916429
$ws.Range("G9").Value = 0.25
$ws.Range("H9").Value = 0.1371845217075572
$ws.Range("I9").Value = "Dunn"
$ws.Range("J9").Value = $false
$ws.Range("K9").Value = -11.5884259259259

# Row 10
$ws.Range("D10").Value = "mannwhitneyu"
$ws.Range("E10").Value = 216
$ws.Range("F10").Value = 0.04854570503465683
$ws.Range("G10").Value = 0.3333333333333334
$ws.Range("H10").Value = 0.04750893270660872
$ws.Range("I10").Value = "Dunn"
$ws.Range("J10").Value = $true
$ws.Range("K10").Value = -4.448708333333332

# Row 11
$ws.Range("D11").Value = "mannwhitneyu"
$ws.Range("E11").Value = 76
$ws.Range("F11").Value = [double]"5.58618703577025e-06"
$ws.Range("G11").Value = 0.7654320987654322
$ws.Range("H11").Value = [double]"5.348038683522629e-06"
$ws.Range("I11").Value = "Dunn"
$ws.Range("J11").Value = $true
$ws.Range("K11").Value = -26.00312037037034

# Row 12
$ws.Range("D12").Value = "mannwhitneyu"
$ws.Range("E12").Value = 167
$ws.Range("F12").Value = 0.004082955523083316
$ws.Range("G12").Value = 0.4845679012345679
$ws.Range("H12").Value = 0.003965978142992611
$ws.Range("I12").Value = "Dunn"
$ws.Range("J12").Value = $true
$ws.Range("K12").Value = -2747.835587962965

# Row 16
$ws.Range("D16").Value = "mannwhitneyu"
$ws.Range("E16").Value = 309
$ws.Range("F16").Value = 0.7901842868705102
$ws.Range("G16").Value = 0.04629629629629628
$ws.Range("H16").Value = 0.7831271046098796
$ws.Range("I16").Value = "Dunn"
$ws.Range("J16").Value = $false
$ws.Range("K16").Value = -5.87546296296297

# Row 17
$ws.Range("D17").Value = "mannwhitneyu"
$ws.Range("E17").Value = 265.5
$ws.Range("F17").Value = 0.2871994947301779
$ws.Range("G17").Value = 0.1805555555555556
$ws.Range("H17").Value = 0.2830647182506413
$ws.Range("I17").Value = "Dunn"
$ws.Range("J17").Value = $false
$ws.Range("K17").Value = -2.401601851851854

# Row 18
$ws.Range("D18").Value = "mannwhitneyu"
$ws.Range("E18").Value = 338
$ws.Range("F18").Value = 0.8043534959231337
$ws.Range("G18").Value = -0.04320987654320985
$ws.Range("H18").Value = 0.7972624975694652
$ws.Range("I18").Value = "Dunn"
$ws.Range("J18").Value = $false
$ws.Range("K18").Value = 15.09937037037039

# Row 19
$ws.Range("D19").Value = "mannwhitneyu"
$ws.Range("E19").Value = 212
$ws.Range("F19").Value = 0.04076053620183384
$ws.Range("G19").Value = 0.345679012345679
$ws.Range("H19").Value = 0.03986622195092659
$ws.Range("I19").Value = "Dunn"
$ws.Range("J19").Value = $true
$ws.Range("K19").Value = -1186.960046296297

# Row 23
$ws.Range("D23").Value = "mannwhitneyu"
$ws.Range("E23").Value = 542
$ws.Range("F23").Value = [double]"6.579704265068806e-05"
$ws.Range("G23").Value = -0.6728395061728396
$ws.Range("H23").Value = [double]"6.329743910048687e-05"
$ws.Range("I23").Value = "Dunn"
$ws.Range("J23").Value = $true
$ws.Range("K23").Value = 27.66805555555555

# Row 24
$ws.Range("D24").Value = "mannwhitneyu"
$ws.Range("E24").Value = 577
$ws.Range("F24").Value = [double]"3.600172181692865e-06"
$ws.Range("G24").Value = -0.7808641975308641
$ws.Range("H24").Value = [double]"3.443899808455802e-06"
$ws.Range("I24").Value = "Dunn"
$ws.Range("J24").Value = $true
$ws.Range("K24").Value = 9.711703703703702

# Row 25
$ws.Range("D25").Value = "mannwhitneyu"
$ws.Range("E25").Value = 562
$ws.Range("F25").Value = [double]"1.312674343949274e-05"
$ws.Range("G25").Value = -0.7345679012345678
$ws.Range("H25").Value = [double]"1.258746709260664e-05"
$ws.Range("I25").Value = "Dunn"
$ws.Range("J25").Value = $true
$ws.Range("K25").Value = 27.45413425925926

# Row 26
$ws.Range("D26").Value = "mannwhitneyu"
$ws.Range("E26").Value = 610
$ws.Range("F26").Value = [double]"1.616656022567341e-07"
$ws.Range("G26").Value = -0.882716049382716
$ws.Range("H26").Value = [double]"1.538199871942933e-07"
$ws.Range("I26").Value = "Dunn"
$ws.Range("J26").Value = $true
$ws.Range("K26").Value = 4737.355986111112

# Row 30
$ws.Range("D30").Value = "mannwhitneyu"
$ws.Range("E30").Value = 648
$ws.Range("F30").Value = [double]"2.920070438311818e-09"
$ws.Range("G30").Value = -1
$ws.Range("H30").Value = [double]"2.761122110912224e-09"
$ws.Range("I30").Value = "Dunn"
$ws.Range("J30").Value = $true
$ws.Range("K30").Value = 135.0351851851852

# Row 31
$ws.Range("D31").Value = "mannwhitneyu"
$ws.Range("E31").Value = 648
$ws.Range("F31").Value = [double]"2.920070438311818e-09"
$ws.Range("G31").Value = -1
$ws.Range("H31").Value = [double]"2.761122110912224e-09"
$ws.Range("I31").Value = "Dunn"
$ws.Range("J31").Value = $true
$ws.Range("K31").Value = 35.42746296296296

# Row 32
$ws.Range("D32").Value = "mannwhitneyu"
$ws.Range("E32").Value = 310
$ws.Range("F32").Value = 0.8043534959231337
$ws.Range("G32").Value = 0.04320987654320985
$ws.Range("H32").Value = 0.7972624975694652
$ws.Range("I32").Value = "Dunn"
$ws.Range("J32").Value = $false
$ws.Range("K32").Value = -0.2727268518518429

# Row 33
$ws.Range("D33").Value = "mannwhitneyu"
$ws.Range("E33").Value = 648
$ws.Range("F33").Value = [double]"2.920070438311818e-09"
$ws.Range("G33").Value = -1
$ws.Range("H33").Value = [double]"2.761122110912224e-09"
$ws.Range("I33").Value = "Dunn"
$ws.Range("J33").Value = $true
$ws.Range("K33").Value = 12919.49476388889

# Row 37
$ws.Range("D37").Value = "mannwhitneyu"
$ws.Range("E37").Value = 419
$ws.Range("F37").Value = 0.0829149018536659
$ws.Range("G37").Value = -0.2932098765432098
$ws.Range("H37").Value = 0.08129997273561756
$ws.Range("I37").Value = "Dunn"
$ws.Range("J37").Value = $false
$ws.Range("K37").Value = 11.29722222222227

# Row 38
$ws.Range("D38").Value = "mannwhitneyu"
$ws.Range("E38").Value = 554
$ws.Range("F38").Value = [double]"2.5394276752008e-05"
$ws.Range("G38").Value = -0.7098765432098766
$ws.Range("H38").Value = [double]"2.438242870832336e-05"
$ws.Range("I38").Value = "Dunn"
$ws.Range("J38").Value = $true
$ws.Range("K38").Value = 12.0933101851852

# Row 39
$ws.Range("D39").Value = "mannwhitneyu"
$ws.Range("E39").Value = 622
$ws.Range("F39").Value = [double]"4.789711871151461e-08"
$ws.Range("G39").Value = -0.9197530864197532
$ws.Range("H39").Value = [double]"4.548338411176395e-08"
$ws.Range("I39").Value = "Dunn"
$ws.Range("J39").Value = $true
$ws.Range("K39").Value = 68.56191666666666

# Row 40
$ws.Range("D40").Value = "mannwhitneyu"
$ws.Range("E40").Value = 606
$ws.Range("F40").Value = [double]"2.399830096728318e-07"
$ws.Range("G40").Value = -0.8703703703703705
$ws.Range("H40").Value = [double]"2.284857699028651e-07"
$ws.Range("I40").Value = "Dunn"
$ws.Range("J40").Value = $true
$ws.Range("K40").Value = 6741.810004629631
